$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new row of data (row 5) mirroring the existing rows 2-4.
$row = 5

$ws.Cells.Item(4, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item($row, 1).Value = 42587.819594907407

$ws.Cells.Item($row, 2).Value = "Named"

$ws.Cells.Item($row, 3).Value = 12637
$ws.Cells.Item($row, 4).Value = 9491
$ws.Cells.Item($row, 5).Value = 568
$ws.Cells.Item($row, 6).Value = 114
$ws.Cells.Item($row, 7).Value = 36
$ws.Cells.Item($row, 8).Value = 75
$ws.Cells.Item($row, 9).Value = 23
$ws.Cells.Item($row, 10).Value = 3
$ws.Cells.Item($row, 11).Value = 2
$ws.Cells.Item($row, 12).Value = 60
$ws.Cells.Item($row, 13).Value = 40
